# 10.8.20 update with changes to time axis
# Append 26 new daily rows (rows 216-241) to the feed log on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: day, date(serial), feedings
$data = @(
  @(215, 44087, 5),
  @(216, 44088, 5),
  @(217, 44089, 5),
  @(218, 44090, 5),
  @(219, 44091, 5),
  @(220, 44092, 5),
  @(221, 44093, 5),
  @(222, 44094, 5),
  @(223, 44095, 5),
  @(224, 44096, 5),
  @(225, 44097, 6),
  @(226, 44098, 6),
  @(227, 44099, 6),
  @(228, 44100, 5),
  @(229, 44101, 6),
  @(230, 44102, 5),
  @(231, 44103, 6),
  @(232, 44104, 6),
  @(233, 44105, 6),
  @(234, 44106, 5),
  @(235, 44107, 5),
  @(236, 44108, 5),
  @(237, 44109, 5),
  @(238, 44110, 5),
  @(239, 44111, 6),
  @(240, 44112, 6)
)

$startRow = 216
$endRow = $startRow + $data.Count - 1

# Match the existing date formatting (style already used by column B, e.g. B2)
# by copying its format onto the new date cells instead of minting a new
# number format.
$ws.Range("B2").Copy()
$ws.Range("B" + $startRow + ":B" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row = $startRow
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# Mirror the scrolled/selected state recorded in the saved view.
$ws.Range("A" + ($endRow + 1)).Select()
